# Generate Report for Handoff
# Update the localization-status report for b.md: it is now ready for
# handoff (previously it was reported as "Handed back: in sync with en-US").

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$handoffDate = "2016-08-31 06:42:02"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusReady
$overview.Range("F3").Value = $statusReady
$overview.Range("G3").Value = $handoffDate

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReady
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-31 06:41:55"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/76a1a69117a7a1a4e3c017e47f366cae11b1a496/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5fc790e9badde33d4dbe4601434eaaf666bfcbbd/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReady
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $handoffDate
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/76a1a69117a7a1a4e3c017e47f366cae11b1a496/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5fc790e9badde33d4dbe4601434eaaf666bfcbbd/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
